# Updated cryptos list on Wed Apr 17 08:57:40 UTC 2024 with GitHub Actions
# Applies updated Price (column D) and Volume(1h) (column E) values for rows 2-51.
# Some Price values are plain decimal numbers (e.g. "543.72"); Excel would
# normally auto-convert such strings to numeric cells when assigned via
# Range.Value. To keep them as text (matching the original inlineStr cells)
# the cell is temporarily switched to a Text number format before the
# assignment and then reset to the default "Normal" style afterwards so no
# residual formatting is left on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.358.15"
$ws.Range("E2").Value = "  +0.09%  "
$ws.Range("D3").Value = "3.074.73"
$ws.Range("E3").Value = "  -0.44%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "543.72"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.67%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.35"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.18%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("D8").Value = "3.069.85"
$ws.Range("E8").Value = "  -0.30%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.497"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.03%  "
$ws.Range("E10").Value = "  +0.37%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.39"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.33%  "
$ws.Range("E12").Value = "  -3.26%  "
$ws.Range("E13").Value = "  +3.25%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.87"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.61%  "
$ws.Range("D15").Value = "3.573.46"
$ws.Range("E15").Value = "  -0.43%  "
$ws.Range("D16").Value = "63.309.63"
$ws.Range("E16").Value = "  +0.08%  "
$ws.Range("E17").Value = "  +1.04%  "
$ws.Range("D18").Value = "3.075.74"
$ws.Range("E18").Value = "  -0.47%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.63"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.53%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "474.57"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.83%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.43"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.57%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.697"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.76%  "
$ws.Range("E23").Value = "  -2.60%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "78.51"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.53%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.18"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.59%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.06%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.71"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.48%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.92"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -6.36%  "
$ws.Range("E29").Value = "  +0.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "26.17"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.76%  "
$ws.Range("E31").Value = "  -4.28%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.16"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.21%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "59.28"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.64%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.30"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -7.95%  "
$ws.Range("E35").Value = "  +7.19%  "
$ws.Range("E36").Value = "  -0.53%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "487.28"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.65%  "
$ws.Range("D38").Value = "3.261.21"
$ws.Range("E38").Value = "  +3.68%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0402"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.60%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0794"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.98%  "
$ws.Range("E41").Value = "  -1.60%  "
$ws.Range("E42").Value = "  -0.45%  "
$ws.Range("E43").Value = "  -2.12%  "
$ws.Range("E44").Value = "  -2.21%  "
$ws.Range("E45").Value = "  +0.09%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "25.43"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.85%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "123.26"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.96%  "
$ws.Range("E48").Value = "  -2.09%  "
$ws.Range("D49").Value = "0.0₃0528"
$ws.Range("E49").Value = "  +5.06%  "
$ws.Range("E50").Value = "  +0.38%  "
$ws.Range("E51").Value = "  -0.60%  "
